# Adding the changes we made on may 9th
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-3.367526054382324, 0.915987193584442, -1.626443386077881),
    @(-1.118759155273438, 2.792432069778442, -4.963344097137451),
    @(-1.557756900787354, -0.5582034587860107, -0.2619988918304443),
    @(1.715949058532715, -1.576748490333557, 5.096891403198242),
    @(-0.8243503570556641, 0.5943599939346313, 1.927432060241699),
    @(-0.2981023788452148, 1.024843096733093, 0.8517363667488098),
    @(0.4371089935302734, 0.3337190449237823, -0.154114544391632),
    @(-0.5641984939575195, -0.3292053341865539, -0.326197862625122),
    @(0.131052017211914, 0.5107872486114502, 0.0540084838867187),
    @(-0.1777238845825195, 0.4102384448051452, 0.1352127194404602)
)

$startRow = 22
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
